$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns I ("I0") and J ("IF") - header cells get the same bold/bordered
# style (s="1") already used by the other header cells (copy format from H1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# I2:J65 data values (row -> I, J)
$iValues = @{}
$jValues = @{}
$iValues[2] = 6; $jValues[2] = 6
$iValues[3] = 7; $jValues[3] = 7
$iValues[4] = 6; $jValues[4] = 6
$iValues[5] = 9; $jValues[5] = 9
$iValues[6] = 6; $jValues[6] = 6
$iValues[7] = 6; $jValues[7] = 8
$iValues[8] = 7; $jValues[8] = 7
$iValues[9] = 5; $jValues[9] = 6
$iValues[10] = 6; $jValues[10] = 7
$iValues[11] = 6; $jValues[11] = 7
$iValues[12] = 8; $jValues[12] = 8
$iValues[13] = 8; $jValues[13] = 8
$iValues[14] = 7; $jValues[14] = 8
$iValues[15] = 9; $jValues[15] = 9
$iValues[16] = 5; $jValues[16] = 6
$iValues[17] = 7; $jValues[17] = 7
$iValues[18] = 9; $jValues[18] = 9
$iValues[19] = 5; $jValues[19] = 6
$iValues[20] = 6; $jValues[20] = 6
$iValues[21] = 7; $jValues[21] = 7
$iValues[22] = 6; $jValues[22] = 7
$iValues[23] = 6; $jValues[23] = 7
$iValues[24] = 8; $jValues[24] = 8
$iValues[25] = 5; $jValues[25] = 6
$iValues[26] = 8; $jValues[26] = 8
$iValues[27] = 7; $jValues[27] = 8
$iValues[28] = 7; $jValues[28] = 8
$iValues[29] = 9; $jValues[29] = 10
$iValues[30] = 6; $jValues[30] = 6
$iValues[31] = 5; $jValues[31] = 5
$iValues[32] = 5; $jValues[32] = 7
$iValues[33] = 9; $jValues[33] = 9
$iValues[34] = 7; $jValues[34] = 7
$iValues[35] = 6; $jValues[35] = 6
$iValues[36] = 8; $jValues[36] = 9
$iValues[37] = 7; $jValues[37] = 7
$iValues[38] = 5; $jValues[38] = 6
$iValues[39] = 8; $jValues[39] = 9
$iValues[40] = 6; $jValues[40] = 6
$iValues[41] = 7; $jValues[41] = 8
$iValues[42] = 9; $jValues[42] = 9
$iValues[43] = 7; $jValues[43] = 7
$iValues[44] = 7; $jValues[44] = 7
$iValues[45] = 7; $jValues[45] = 8
$iValues[46] = 6; $jValues[46] = 6
$iValues[47] = 5; $jValues[47] = 6
$iValues[48] = 8; $jValues[48] = 8
$iValues[49] = 6; $jValues[49] = 6
$iValues[50] = 9; $jValues[50] = 9
$iValues[51] = 7; $jValues[51] = 7
$iValues[52] = 7; $jValues[52] = 8
$iValues[53] = 7; $jValues[53] = 7
$iValues[54] = 9; $jValues[54] = 9
$iValues[55] = 7; $jValues[55] = 8
$iValues[56] = 8; $jValues[56] = 9
$iValues[57] = 8; $jValues[57] = 8
$iValues[58] = 5; $jValues[58] = 5
$iValues[59] = 7; $jValues[59] = 7
$iValues[60] = 5; $jValues[60] = 5
$iValues[61] = 7; $jValues[61] = 7
$iValues[62] = 7; $jValues[62] = 7
$iValues[63] = 5; $jValues[63] = 5
$iValues[64] = 6; $jValues[64] = 6
$iValues[65] = 8; $jValues[65] = 8

for ($row = 2; $row -le 65; $row++) {
  $ws.Cells.Item($row, 9).Value = $iValues[$row]
  $ws.Cells.Item($row, 10).Value = $jValues[$row]
}
